$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13: CamOffestPos
$ws.Range("A13").Value = "CamOffestPos"
$ws.Range("B13").Value = "string"
$ws.Range("C13").Value = $false
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = $true
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = "Friend"
$ws.Range("J13").Value = "acctorid"

# Row 14: CamOffestRot
$ws.Range("A14").Value = "CamOffestRot"
$ws.Range("B14").Value = "string"
$ws.Range("C14").Value = $false
$ws.Range("D14").Value = $false
$ws.Range("E14").Value = $false
$ws.Range("F14").Value = $true
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "Friend"
$ws.Range("J14").Value = "acctorid"

# Copy style from row 12 cells to row 13/14 matching cells (A,B,I,J use style 1)
$ws.Range("A12:J12").Copy()
$ws.Range("A13:J13").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A12:J12").Copy()
$ws.Range("A14:J14").PasteSpecial(-4122) # xlPasteFormats

# The existing list validation on column F ("TRUE,FALSE") covered F13/F14
# already (as part of "F13:F1048576"). Re-apply it explicitly over F13:F14
# so the new rows carry their own validation entry, same as every other
# populated row, and the remaining range is split around them.
$ws.Range("F13:F14").Validation.Delete()
$ws.Range("F13:F14").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Update selection to A14
$ws.Range("A14").Select()
